# Insert a new column B ("Variável") into the worksheet, shifting the
# existing "Valor" and "Colocação" columns to C and D respectively, and
# fill the new column with "Diferença 2022-2013" for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before the current column B (shifts B->C, C->D)
$ws.Columns("B:B").Insert()

# Header for the new column
$ws.Range("B1").Value = "Variável"

# Fill the new column's data rows (2 through 10) with the constant label
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 2).Value = "Diferença 2022-2013"
}
